$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.423.67'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.08%  '

$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.854.44'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.44%  '

$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.33%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.31%  '

$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6316'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.80%  '

$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.11%  '

$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07566'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.19%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2924'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.31%  '

$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.58'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.93%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07756'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.20%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.853.03'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.12%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.035'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.18%  '

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6843'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.64%  '

$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001045'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.57%  '

$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.37'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.04%  '

$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.142'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.39%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.419.50'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.01%  '

$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '230.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.91%  '

$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.09%  '

$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.15%  '

$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.511'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.40%  '

$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.03%  '

$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '159.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.85%  '

$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1398'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.95%  '

$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.459'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.77%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.73'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.26%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.413'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.37%  '

$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.482'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.39%  '

$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05689'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.58%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.152'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.89%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.058'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.82%  '

$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.827'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.70%  '

$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.157'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.33%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6990'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.25%  '

$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.589'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.11%  '

$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.250.39'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.28%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01827'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.30%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.771'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.44%  '

$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.520'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.39%  '

$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9092'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.62%  '

$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.000'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.02%  '

$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.018.03'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.71%  '

$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.55'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.16%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.11'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.31%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.142'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.06%  '

$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1166'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.13%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.067'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.74%  '

$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000116'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.34%  '

$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3968'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.06%  '

$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.675'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.28%  '
